# Insert a new row before row 4, shifting the "realises" (BFO:0000055) row
# down to row 5, then populate the new row 4 with the "informal label"
# annotation property.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(4).Insert()

$ws.Cells.Item(4, 1).Value = ""
$ws.Cells.Item(4, 2).Value = "informal label"
$ws.Cells.Item(4, 3).Value = ""
$ws.Cells.Item(4, 4).Value = ""
$ws.Cells.Item(4, 5).Value = ""
$ws.Cells.Item(4, 6).Value = ""
$ws.Cells.Item(4, 7).Value = ""
$ws.Cells.Item(4, 8).Value = "AnnotationProperty"
